# Auto-generated edit script: update crypto price/volume data per commit
# "Updated cryptos list on Thu Aug 15 13:53:59 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.076.42"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "2.638.51"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'527.56"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").Value = "'144.56"
$ws.Range("E6").Value = "  -1.87%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'0.568"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").Value = "'6.64"
$ws.Range("E9").Value = "  -4.21%  "
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "3.106.86"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").Value = "59.034.18"
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.752.63"
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'20.99"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000137"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "'341.57"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "'4.44"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").Value = "'10.54"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").Value = "'6.33"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'65.14"
$ws.Range("E23").Value = "  +2.55%  "
$ws.Range("D24").Value = "'0.417"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").Value = "'0.167"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "0.0₃0798"
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("D29").Value = "'6.45"
$ws.Range("E29").Value = "  -4.37%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("D32").Value = "'18.92"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").Value = "'150.05"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "'4.20"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").Value = "'1.20"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").Value = "'0.926"
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("D37").Value = "'0.872"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("D39").Value = "'36.52"
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("D41").Value = "'0.997"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").Value = "'0.603"
$ws.Range("E42").Value = "  -5.18%  "
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").Value = "'271.50"
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("D45").Value = "'19.36"
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("D46").Value = "'0.0538"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("D48").Value = "2.048.62"
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D49").Value = "'4.74"
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").Value = "'18.93"
$ws.Range("E51").Value = "  -0.70%  "
